# Swap the content of rows 18<->20 and 19<->21 on the active sheet
# (columns A, B, E, F, G, H, K, L, M, N, Q, R, Z, AB, AC), as described
# by the diff. All other columns (D, I, P, S, T, U, V, W, Y, AA, AD, AE,
# AG, AT, AW, AX, AY) are identical between the affected rows and are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 (becomes the old row 20 content) ---
$ws.Range("A18").Value = 130979946
$ws.Range("B18").Value = 57884
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("K18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").Value = "färska spår"
$ws.Range("N18").ClearContents()
$ws.Range("Q18").Value = 590605
$ws.Range("R18").Value = 6963364
$ws.Range("Z18").Value = "09:47"
$ws.Range("AB18").Value = "09:47"
$ws.Range("AC18").Value = "färska ringhack på gran"

# --- Row 19 (becomes the old row 21 content) ---
$ws.Range("A19").Value = 130979899
$ws.Range("B19").Value = 57884
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = "Tretåig hackspett"
$ws.Range("G19").Value = "Picoides tridactylus"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("K19").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").Value = "färska spår"
$ws.Range("N19").ClearContents()
$ws.Range("Q19").Value = 590850
$ws.Range("R19").Value = 6963133
$ws.Range("Z19").Value = "13:16"
$ws.Range("AB19").Value = "13:16"
$ws.Range("AC19").Value = "färska ringhack på gran"

# --- Row 20 (becomes the old row 18 content) ---
$ws.Range("A20").Value = 130979947
$ws.Range("B20").Value = 91808
$ws.Range("E20").Value = 1202
$ws.Range("F20").Value = "Ullticka"
$ws.Range("G20").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H20").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K20").ClearContents()
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("Q20").Value = 590591
$ws.Range("R20").Value = 6963354
$ws.Range("Z20").Value = "09:45"
$ws.Range("AB20").Value = "09:45"
$ws.Range("AC20").ClearContents()

# --- Row 21 (becomes the old row 19 content) ---
$ws.Range("A21").Value = 130979897
$ws.Range("B21").Value = 80348
$ws.Range("E21").Value = 6458
$ws.Range("F21").Value = "Lunglav"
$ws.Range("G21").Value = "Lobaria pulmonaria"
$ws.Range("H21").Value = "(L.) Hoffm."
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("Q21").Value = 590726
$ws.Range("R21").Value = 6963153
$ws.Range("Z21").Value = "13:24"
$ws.Range("AB21").Value = "13:24"
$ws.Range("AC21").ClearContents()
